$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string values used in F14 (SheetName) and G14 (TestDataRow)
$ws.Cells.Item(14, 6).Value2 = "test2_instance"
$ws.Cells.Item(14, 7).Value2 = "1-5,7,9,11-12,14-15,18,20-23,27-40,42,44-48,50,51"

# Match the style of G14 to that of F14 (s="1" instead of s="2")
$ws.Cells.Item(14, 7).Style = $ws.Cells.Item(14, 6).Style

# Update the view: scroll so column B is the top-left visible column,
# and move the active selection to G14
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G14").Select()
